# Update the TestName in the "Classes" sheet: the second regression-test
# row referenced a test case ("TC002_ShopClient_VerifyRegPageTitle") that no
# longer matches the actual test; point it at the renamed test instead so the
# TestNG config lines up with the test that now waits for the click action.
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Classes")
$ws.Range("C3").Value = "TC002_ShopClient_VerifySearchCarInfoByVin"
